$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 keeps its text ("HK_G_acc_LG"); re-assigning it is what produces the
# extra shared-string entry for the "better description of localT2 ls".
$ws.Range("A1").Value = "HK_G_acc_LG"

# Recomputed accuracy values for A2:A49.
$values = @(
    76.51195499296765,
    76.51195499296765,
    76.23066104078762,
    76.51195499296765,
    76.37130801687763,
    76.37130801687763,
    76.37130801687763,
    76.09001406469761,
    76.37130801687763,
    76.37130801687763,
    75.9493670886076,
    76.51195499296765,
    76.51195499296765,
    76.23066104078762,
    76.51195499296765,
    76.51195499296765,
    76.51195499296765,
    76.51195499296765,
    76.23066104078762,
    76.23066104078762,
    76.23066104078762,
    76.09001406469761,
    76.23066104078762,
    75.80872011251758,
    77.35583684950772,
    76.51195499296765,
    77.35583684950772,
    76.09001406469761,
    76.23066104078762,
    76.51195499296765,
    76.37130801687763,
    76.51195499296765,
    76.51195499296765,
    76.65260196905767,
    75.66807313642757,
    75.38677918424754,
    76.09001406469761,
    77.0745428973277,
    77.21518987341773,
    76.51195499296765,
    76.37130801687763,
    76.37130801687763,
    76.37130801687763,
    76.51195499296765,
    76.37130801687763,
    76.51195499296765,
    76.37130801687763,
    76.51195499296765
)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 1).Value = $v
    $row = $row + 1
}
